$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column (H1), copying formatting (style) from G1 header cell
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Save values for rows 2-10
$saveValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 0
    6 = 1
    7 = 0
    8 = 1
    9 = 0
    10 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}
